$d = $word.ActiveDocument

$ids = @("p066v_4", "p067r_1", "p067r_2", "p067r_3", "p067r_4")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"
    $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}
